# Update the "Förändrad" (Changed) date column (column C) for all data rows
# from serial date 45775 (2025-04-28) to 45776 (2025-04-29).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 43; $row++) {
    $ws.Cells.Item($row, 3).Value = 45776
}
